$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.324917666666666
$ws.Range("H2").Value = 24.974753
$ws.Range("I2").Value = 0.8193616330571973
$ws.Range("J2").Value = 0.8193616330571972
$ws.Range("M2").Value = 0.9705896666666667
$ws.Range("N2").Value = 2.911769
$ws.Range("O2").Value = 0.02073452941466921
$ws.Range("P2").Value = 0.02073452941466921
$ws.Range("Q2").Value = 8.080079063117443
$ws.Range("R2").Value = 72.720711568057
$ws.Range("S2").Value = 0.01698907788187586
$ws.Range("T2").Value = 0.01698907788187586

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.324917666666666
$ws.Range("H3").Value = 24.974753
$ws.Range("I3").Value = 0.8193616330571973
$ws.Range("J3").Value = 0.8193616330571972
$ws.Range("O3").Value = 0.5628689972673966
$ws.Range("P3").Value = 0.5628689972673966
$ws.Range("Q3").Value = 219.3455134255699
$ws.Range("R3").Value = 1974.10962083013
$ws.Range("S3").Value = 0.4611932607982812
$ws.Range("T3").Value = 0.4611932607982812

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.324917666666666
$ws.Range("H4").Value = 24.974753
$ws.Range("I4").Value = 0.8193616330571973
$ws.Range("J4").Value = 0.8193616330571972
$ws.Range("M4").Value = 19.49164633333333
$ws.Range("N4").Value = 58.47493899999999
$ws.Range("O4").Value = 0.4163964733179342
$ws.Range("P4").Value = 0.4163964733179341
$ws.Range("Q4").Value = 162.2663509127852
$ws.Range("R4").Value = 1460.397158215067
$ws.Range("S4").Value = 0.3411792943770403
$ws.Range("T4").Value = 0.3411792943770401

# Row 5
$ws.Range("G5").Value = 1.835330666666667
$ws.Range("H5").Value = 5.505992
$ws.Range("I5").Value = 0.1806383669428028
$ws.Range("J5").Value = 0.1806383669428027
$ws.Range("M5").Value = 0.9705896666666667
$ws.Range("N5").Value = 2.911769
$ws.Range("O5").Value = 0.02073452941466921
$ws.Range("P5").Value = 0.02073452941466921
$ws.Range("Q5").Value = 1.781352979983111
$ws.Range("R5").Value = 16.032176819848
$ws.Range("S5").Value = 0.003745451532793355
$ws.Range("T5").Value = 0.003745451532793354

# Row 6
$ws.Range("G6").Value = 1.835330666666667
$ws.Range("H6").Value = 5.505992
$ws.Range("I6").Value = 0.1806383669428028
$ws.Range("J6").Value = 0.1806383669428027
$ws.Range("O6").Value = 0.5628689972673966
$ws.Range("P6").Value = 0.5628689972673966
$ws.Range("Q6").Value = 48.35742087847999
$ws.Range("R6").Value = 435.21678790632
$ws.Range("S6").Value = 0.1016757364691154
$ws.Range("T6").Value = 0.1016757364691154

# Row 7
$ws.Range("G7").Value = 1.835330666666667
$ws.Range("H7").Value = 5.505992
$ws.Range("I7").Value = 0.1806383669428028
$ws.Range("J7").Value = 0.1806383669428027
$ws.Range("M7").Value = 19.49164633333333
$ws.Range("N7").Value = 58.47493899999999
$ws.Range("O7").Value = 0.4163964733179342
$ws.Range("P7").Value = 0.4163964733179341
$ws.Range("Q7").Value = 35.77361625938755
$ws.Range("R7").Value = 321.962546334488
$ws.Range("S7").Value = 0.07521717894089398
$ws.Range("T7").Value = 0.07521717894089396
